$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1909650924024641
$ws.Range("C2").Value = 0.5585215605749486
$ws.Range("J2").Value = 0.01232032854209446
$ws.Range("P2").Value = 0.1314168377823409
$ws.Range("S2").Value = 0.106776180698152
$ws.Range("B3").Value = 0.01071428571428571
$ws.Range("C3").Value = 0.02142857142857143
$ws.Range("J3").Value = 0.01785714285714286
$ws.Range("P3").Value = 0.7392857142857143
$ws.Range("S3").Value = 0.2107142857142857
$ws.Range("J4").Value = 0.03061224489795918
$ws.Range("P4").Value = 0.6836734693877551
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.05527638190954774
$ws.Range("D6").Value = 0.01758793969849246
$ws.Range("F6").Value = 0.0678391959798995
$ws.Range("J6").Value = 0.2462311557788945
$ws.Range("O6").Value = 0.01256281407035176
$ws.Range("Q6").Value = 0.1758793969849246
$ws.Range("R6").Value = 0.07537688442211055
$ws.Range("S6").Value = 0.3492462311557789
$ws.Range("B7").Value = 0.09832134292565947
$ws.Range("D7").Value = 0.026378896882494
$ws.Range("E7").Value = 0.004796163069544364
$ws.Range("F7").Value = 0.05035971223021583
$ws.Range("J7").Value = 0.1270983213429256
$ws.Range("O7").Value = 0.01199040767386091
$ws.Range("Q7").Value = 0.2206235011990408
$ws.Range("R7").Value = 0.03836930455635491
$ws.Range("S7").Value = 0.4220623501199041
$ws.Range("B8").Value = 0.07427937915742794
$ws.Range("D8").Value = 0.02549889135254989
$ws.Range("F8").Value = 0.06430155210643015
$ws.Range("J8").Value = 0.1130820399113082
$ws.Range("O8").Value = 0.0188470066518847
$ws.Range("Q8").Value = 0.2028824833702882
$ws.Range("R8").Value = 0.07538802660753881
$ws.Range("S8").Value = 0.4257206208425721
$ws.Range("B9").Value = 0.08994708994708994
$ws.Range("D9").Value = 0.02116402116402116
$ws.Range("F9").Value = 0.05291005291005291
$ws.Range("J9").Value = 0.1084656084656085
$ws.Range("O9").Value = 0.02116402116402116
$ws.Range("Q9").Value = 0.1719576719576719
$ws.Range("R9").Value = 0.08994708994708994
$ws.Range("S9").Value = 0.4444444444444444
$ws.Range("B10").Value = 0.09405728943993159
$ws.Range("D10").Value = 0.02265925609234716
$ws.Range("E10").Value = 0.0008550662676357417
$ws.Range("F10").Value = 0.06584010260795212
$ws.Range("J10").Value = 0.1120136810602822
$ws.Range("O10").Value = 0.013253527148354
$ws.Range("Q10").Value = 0.2248824283882001
$ws.Range("R10").Value = 0.07695596408721676
$ws.Range("S10").Value = 0.3894826849080804
$ws.Range("G11").Value = 0.1261398176291793
$ws.Range("J11").Value = 0.08966565349544073
$ws.Range("K11").Value = 0.1899696048632219
$ws.Range("L11").Value = 0.5851063829787234
$ws.Range("S11").Value = 0.00911854103343465
$ws.Range("G12").Value = 0.7398989898989899
$ws.Range("J12").Value = 0.1868686868686869
$ws.Range("K12").Value = 0.005050505050505051
$ws.Range("L12").Value = 0.02272727272727273
$ws.Range("S12").Value = 0.04545454545454546
$ws.Range("F15").Value = 0.01351351351351351
$ws.Range("H15").Value = 0.1689189189189189
$ws.Range("I15").Value = 0.08558558558558559
$ws.Range("J15").Value = 0.3198198198198198
$ws.Range("K15").Value = 0.08783783783783784
$ws.Range("M15").Value = 0.006756756756756757
$ws.Range("O15").Value = 0.0945945945945946
$ws.Range("S15").Value = 0.222972972972973
$ws.Range("F16").Value = 0.01834862385321101
$ws.Range("H16").Value = 0.1743119266055046
$ws.Range("I16").Value = 0.09174311926605505
$ws.Range("J16").Value = 0.3914373088685015
$ws.Range("K16").Value = 0.1345565749235474
$ws.Range("M16").Value = 0.01223241590214067
$ws.Range("N16").Value = 0.003058103975535168
$ws.Range("O16").Value = 0.06116207951070336
$ws.Range("S16").Value = 0.1131498470948012
$ws.Range("F17").Value = 0.01929260450160772
$ws.Range("H17").Value = 0.1939978563772776
$ws.Range("I17").Value = 0.09646302250803858
$ws.Range("J17").Value = 0.3762057877813505
$ws.Range("K17").Value = 0.1018220793140407
$ws.Range("M17").Value = 0.01714898177920686
$ws.Range("N17").Value = 0.001071811361200429
$ws.Range("O17").Value = 0.05894962486602358
$ws.Range("S17").Value = 0.135048231511254
$ws.Range("F18").Value = 0.01829268292682927
$ws.Range("H18").Value = 0.1676829268292683
$ws.Range("I18").Value = 0.08841463414634146
$ws.Range("J18").Value = 0.4298780487804878
$ws.Range("K18").Value = 0.0975609756097561
$ws.Range("M18").Value = 0.009146341463414634
$ws.Range("O18").Value = 0.04878048780487805
$ws.Range("S18").Value = 0.1402439024390244
$ws.Range("F19").Value = 0.01636151149201403
$ws.Range("H19").Value = 0.2095831710167511
$ws.Range("I19").Value = 0.07635372029606545
$ws.Range("J19").Value = 0.3393065835605765
$ws.Range("K19").Value = 0.1238800155823919
$ws.Range("M19").Value = 0.02726918582002337
$ws.Range("N19").Value = 0.001558239189715621
$ws.Range("O19").Value = 0.07128944292948967
$ws.Range("S19").Value = 0.1343981301129723
